# ---------------------------------------------------------------------------
# Adds a new "2022-Q3" quarter to the 689009-九号公司 workbook:
#   1. "总计" (summary) sheet gets a new row 2 (2022-Q3 / 26 / 26.56) and all
#      the previously existing quarters shift down by one row.
#   2. A brand-new worksheet named "2022-Q3" is inserted right after "总计"
#      (pushing every other quarter tab one slot to the right) and is filled
#      with the fund-holding breakdown for that quarter.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert the 2022-Q3 row at the
#    top of the data block and shift the rest down by one row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Snapshot the existing data rows (rows 2-8) before we overwrite them.
$existingRows = New-Object System.Collections.ArrayList
for ($r = 2; $r -le 8; $r++) {
    [void]$existingRows.Add(@(
        $summary.Cells.Item($r, 2).Value(),
        $summary.Cells.Item($r, 3).Value(),
        $summary.Cells.Item($r, 4).Value()
    ))
}

# New row order: 2022-Q3 first, then the previously existing rows.
$newRows = New-Object System.Collections.ArrayList
[void]$newRows.Add(@("2022-Q3", 26, 26.56))
foreach ($row in $existingRows) {
    [void]$newRows.Add($row)
}

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $i + 2
    $row = $newRows[$i]
    $summary.Cells.Item($r, 1).Value = $i
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
}

# Row 9 is brand new - copy the bold/bordered index-column style from A8.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计".
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Header row (row 1). Columns B:H, styled like the other quarter sheets.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, 2 + $c).Value = $headers[$c]
}

# Fund-holding rows (row 2 .. 27).
$funds = @(
    @("011058", "景顺长城成长龙头一年持有期混合A", "51.58", "91.71", "6.08", "3.1361", 3),
    @("011329", "景顺长城新能源产业股票C",         "46.25", "89.51", "6.24", "2.8860", 2),
    @("260101", "景顺长城优选混合",                 "49.69", "76.36", "5.62", "2.7926", 1),
    @("001975", "景顺长城环保优势股票",             "44.53", "92.46", "6.22", "2.7698", 1),
    @("007802", "兴全合泰混合A",                   "66.03", "90.33", "3.61", "2.3837", 8),
    @("006435", "景顺长城创新成长混合",             "36.71", "92.71", "5.93", "2.1769", 1),
    @("014639", "兴证全球合衡三年持有混合A",       "53.34", "81.42", "3.57", "1.9042", 6),
    @("011328", "景顺长城新能源产业股票A",         "26.50", "89.51", "6.24", "1.6536", 2),
    @("009795", "嘉实远见精选两年持有期混合",       "55.77", "93.56", "2.88", "1.6062", 9),
    @("000970", "东方红睿元三年定期开放灵活配置混合", "31.77", "72.56", "3.31", "1.0516", 7),
    @("506006", "汇添富科创板2年定期开放混合",     "16.61", "88.62", "5.36", "0.8903", 5),
    @("009376", "景顺长城成长领航混合",             "13.15", "92.87", "5.90", "0.7758", 1),
    @("007803", "兴全合泰混合C",                   "12.16", "90.33", "3.61", "0.4390", 8),
    @("011059", "景顺长城成长龙头一年持有期混合C", "6.31",  "91.71", "6.08", "0.3836", 3),
    @("009683", "汇添富创新增长一年定期开放混合A", "9.83",  "80.10", "3.82", "0.3755", 4),
    @("000480", "东方红新动力灵活配置混合",         "12.63", "78.30", "2.94", "0.3713", 9),
    @("001564", "东方红京东大数据灵活配置混合",     "8.84",  "73.95", "3.07", "0.2714", 9),
    @("004476", "景顺长城沪港深领先科技股票",       "9.71",  "81.72", "2.39", "0.2321", 9),
    @("260111", "景顺长城公司治理混合",             "3.52",  "91.70", "5.54", "0.1950", 1),
    @("014640", "兴证全球合衡三年持有混合C",       "3.69",  "81.42", "3.57", "0.1317", 6),
    @("009684", "汇添富创新增长一年定期开放混合C", "1.19",  "80.10", "3.82", "0.0455", 4),
    @("562500", "华夏中证机器人ETF",               "1.57",  "99.51", "2.24", "0.0352", 9),
    @("159770", "天弘中证机器人ETF",               "1.04",  "99.68", "2.24", "0.0233", 9),
    @("562360", "银华中证机器人ETF",               "0.77",  "97.23", "2.20", "0.0169", 9),
    @("001535", "景顺长城改革机遇灵活配置混合A",   "0.26",  "64.66", "4.35", "0.0113", 1),
    @("007945", "景顺长城改革机遇灵活配置混合C",   "0.06",  "64.66", "4.35", "0.0026", 1)
)

$lastRow = 1 + $funds.Length

# Columns B (fund code), D, E, F, G hold values that look numeric (and some
# have leading/trailing zeros, e.g. "011058" or "26.50") but must stay text,
# exactly like the sibling quarter sheets. Force text formatting up front so
# the values below aren't silently coerced into numbers.
$q3.Range("B2:B$lastRow").NumberFormat = "@"
$q3.Range("D2:G$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $funds.Length; $i++) {
    $r = $i + 2
    $fund = $funds[$i]
    $q3.Cells.Item($r, 1).Value = $i          # A: running index (number)
    $q3.Cells.Item($r, 2).Value = $fund[0]    # B: fund code (text)
    $q3.Cells.Item($r, 3).Value = $fund[1]    # C: fund name (text)
    $q3.Cells.Item($r, 4).Value = $fund[2]    # D: fund scale (text)
    $q3.Cells.Item($r, 5).Value = $fund[3]    # E: stock position (text)
    $q3.Cells.Item($r, 6).Value = $fund[4]    # F: position ratio (text)
    $q3.Cells.Item($r, 7).Value = $fund[5]    # G: holding value (text)
    $q3.Cells.Item($r, 8).Value = $fund[6]    # H: position rank (number)
}

# Match the bold/bordered style used on the header row + index column of the
# sibling quarter sheets (e.g. "2022-Q2"), copied straight from "总计".
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$summary.Range("A2").Copy()
$q3.Range("A2:A$lastRow").PasteSpecial(-4122)

$q3.Range("A1").Select()
